# Locate the unique anchor "10:34" near the end of the document to scope
# subsequent searches (the phrase "documentation technique" also appears
# earlier in the document, so later searches must stay scoped after this
# anchor to avoid ambiguous matches).
$d = $word.ActiveDocument
$anchor = $d.Content
$anchor.Find.Execute("10:34 : Fin de la fonctionnalit", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorStart = $anchor.Start

# Step 1: split the final paragraph in two, right after "cela fonctionne."
# so the trailing bookmark (_GoBack) ends up alone in its own paragraph.
$splitScope = $d.Range($anchorStart, $d.Content.End)
$splitScope.Find.Execute("cela fonctionne.", $true, $false, $false, $false, $false, $true, 1, $false, "cela fonctionne.`r", 2)

# Step 2: append the new sentence as a new run right after "cela fonctionne."
# (now the first sentence of the newly split first paragraph).
$insertScope = $d.Range($anchorStart, $d.Content.End)
$insertScope.Find.Execute("cela fonctionne.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertScope.Collapse(0)
$insertScope.Bold = 1
$insertScope.InsertAfter(" J’avance la documentation technique.")
$insertScope.Bold = 0
